# Auto-generated Excel COM-interop script to apply the Shinryu_Profits profession-sheet data refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3273.0908
$ws.Range("I76").Value = 2876
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 2876
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -2561
$ws.Range("N76").Value = -4130

$ws.Range("H79").Value = 3273.0908
$ws.Range("I79").Value = 2876
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 2876
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -1784
$ws.Range("N79").Value = -5684

$ws.Range("H111").Value = 58824030
$ws.Range("I111").Value = 66667024
$ws.Range("J111").Value = 1600
$ws.Range("K111").Value = 200001072
$ws.Range("L111").Value = 4800
$ws.Range("M111").Value = -199998005
$ws.Range("N111").Value = -10934

$ws.Range("H129").Value = 716.99
$ws.Range("J129").Value = 716.99
$ws.Range("L129").Value = 2150.97
$ws.Range("N129").Value = -12150.97

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13169218
$ws.Range("I32").Value = 18185318
$ws.Range("J32").Value = 31812.285
$ws.Range("K32").Value = 18185318
$ws.Range("L32").Value = 31812.285
$ws.Range("M32").Value = -18185031
$ws.Range("N32").Value = -32386.285

$ws.Range("H61").Value = 3204.8667
$ws.Range("I61").Value = 2176
$ws.Range("J61").Value = 5262.6
$ws.Range("K61").Value = 2176
$ws.Range("L61").Value = 5262.6
$ws.Range("M61").Value = -1964
$ws.Range("N61").Value = -5686.6

$ws.Range("H74").Value = 4796.5454
$ws.Range("I74").Value = 5634.0435
$ws.Range("J74").Value = 2870.3
$ws.Range("K74").Value = 5634.0435
$ws.Range("L74").Value = 2870.3
$ws.Range("M74").Value = -4760.0435
$ws.Range("N74").Value = -4618.3

$ws.Range("H77").Value = 4796.5454
$ws.Range("I77").Value = 5634.0435
$ws.Range("J77").Value = 2870.3
$ws.Range("K77").Value = 28170.2175
$ws.Range("L77").Value = 14351.5
$ws.Range("M77").Value = -23802.2175
$ws.Range("N77").Value = -23087.5

$ws.Range("H132").Value = 1629.3256
$ws.Range("I132").Value = 1075.1936
$ws.Range("J132").Value = 3060.8333
$ws.Range("K132").Value = 3225.5808
$ws.Range("L132").Value = 9182.499899999999
$ws.Range("M132").Value = -695.5808000000002
$ws.Range("N132").Value = -14242.4999

$ws.Range("H136").Value = 3204.8667
$ws.Range("I136").Value = 2176
$ws.Range("J136").Value = 5262.6
$ws.Range("K136").Value = 6528
$ws.Range("L136").Value = 15787.8
$ws.Range("M136").Value = -3978
$ws.Range("N136").Value = -20887.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3108.4666
$ws.Range("I20").Value = 2471.3333
$ws.Range("J20").Value = 3533.2222
$ws.Range("K20").Value = 2471.3333
$ws.Range("L20").Value = 3533.2222
$ws.Range("M20").Value = -2224.3333
$ws.Range("N20").Value = -4027.2222

$ws.Range("H86").Value = 1780.9524
$ws.Range("I86").Value = 1768.3846
$ws.Range("J86").Value = 1801.375
$ws.Range("K86").Value = 1768.3846
$ws.Range("L86").Value = 1801.375
$ws.Range("M86").Value = -645.3846000000001
$ws.Range("N86").Value = -4047.375

$ws.Range("H89").Value = 1780.9524
$ws.Range("I89").Value = 1768.3846
$ws.Range("J89").Value = 1801.375
$ws.Range("K89").Value = 8841.923000000001
$ws.Range("L89").Value = 9006.875
$ws.Range("M89").Value = -3225.923000000001
$ws.Range("N89").Value = -20238.875

$ws.Range("H94").Value = 682.5
$ws.Range("I94").Value = 628.125
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 628.125
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -177.125
$ws.Range("N94").Value = -1802

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3034.3157
$ws.Range("I62").Value = 2900
$ws.Range("J62").Value = 3183.5557
$ws.Range("K62").Value = 2900
$ws.Range("L62").Value = 3183.5557
$ws.Range("M62").Value = -2276
$ws.Range("N62").Value = -4431.5557

$ws.Range("H65").Value = 3034.3157
$ws.Range("I65").Value = 2900
$ws.Range("J65").Value = 3183.5557
$ws.Range("K65").Value = 14500
$ws.Range("L65").Value = 15917.7785
$ws.Range("M65").Value = -11380
$ws.Range("N65").Value = -22157.7785

$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -7246

$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 25000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -36232

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 27966.072
$ws.Range("I87").Value = 3441.7693
$ws.Range("J87").Value = 31807.229
$ws.Range("K87").Value = 10325.3079
$ws.Range("L87").Value = 95421.68700000001
$ws.Range("M87").Value = -9077.3079
$ws.Range("N87").Value = -97917.68700000001

$ws.Range("H90").Value = 27966.072
$ws.Range("I90").Value = 3441.7693
$ws.Range("J90").Value = 31807.229
$ws.Range("K90").Value = 30975.9237
$ws.Range("L90").Value = 286265.061
$ws.Range("M90").Value = -24735.9237
$ws.Range("N90").Value = -298745.061

$ws.Range("H113").Value = 1567777.2
$ws.Range("I113").Value = 2873911.8
$ws.Range("J113").Value = 416
$ws.Range("K113").Value = 8621735.399999999
$ws.Range("L113").Value = 1248
$ws.Range("M113").Value = -8619565.399999999
$ws.Range("N113").Value = -5588

$ws.Range("H131").Value = 652.4286
$ws.Range("I131").Value = 652.4286
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 1957.2858
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 3082.7142
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 60767.5
$ws.Range("J119").Value = 60767.5
$ws.Range("L119").Value = 60767.5
$ws.Range("N119").Value = -70443.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H112").Value = 44993.5
$ws.Range("J112").Value = 44993.5
$ws.Range("L112").Value = 44993.5
$ws.Range("N112").Value = -47947.5

$ws.Range("H114").Value = 30000
$ws.Range("J114").Value = 30000
$ws.Range("L114").Value = 30000
$ws.Range("N114").Value = -38678

$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
